$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows of questions starting at row 26 (rows 23-25 stay blank, matching
# the existing gap pattern already present in the sheet). Values are entered in
# the same order the author originally typed them (matching the shared-string
# table insertion order), with "GAMS Engine credentials?" entered before the
# last three rows even though it ends up placed at the bottom of the list.
$ws.Range("A26").Value = "Questions for Amit"
$ws.Range("A27").Value = "FLO_EMISS?"
$ws.Range("A28").Value = "COM_PROJ?"
$ws.Range("A32").Value = "GAMS Engine credentials?"
$ws.Range("A29").Value = "VO setup - linking to Git."
$ws.Range("A30").Value = "Is it possible to initiate the VO run from the command prompt?"
$ws.Range("A31").Value = "Does the user get access to the GDX file when run is done?"

# Update the selection / top-left cell to match the new active selection area.
$ws.Activate()
$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 4
